$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the diff (crypto price/volume refresh,
# plus a few rows whose coin/link/price/volume were swapped with neighboring rows).
# Cells that hold numeric-looking text (e.g. "613.52") are written with a leading
# apostrophe so Excel keeps them as text, matching the source workbook's formatting.

$ws.Cells.Item(2, 4).Formula = '69.871.25'
$ws.Cells.Item(2, 5).Formula = '  +3.68%  '
$ws.Cells.Item(3, 4).Formula = '3.776.56'
$ws.Cells.Item(3, 5).Formula = '  +21.12%  '
$ws.Cells.Item(4, 5).Formula = '  +0.13%  '
$ws.Cells.Item(5, 4).Formula = '''613.52'
$ws.Cells.Item(5, 5).Formula = '  +6.43%  '
$ws.Cells.Item(6, 4).Formula = '''176.60'
$ws.Cells.Item(6, 5).Formula = '  -1.39%  '
$ws.Cells.Item(7, 4).Formula = '3.773.91'
$ws.Cells.Item(7, 5).Formula = '  +21.02%  '
$ws.Cells.Item(8, 5).Formula = '  +0.15%  '
$ws.Cells.Item(9, 4).Formula = '''0.543'
$ws.Cells.Item(9, 5).Formula = '  +5.21%  '
$ws.Cells.Item(10, 4).Formula = '''0.166'
$ws.Cells.Item(10, 5).Formula = '  +9.13%  '
$ws.Cells.Item(11, 4).Formula = '''6.36'
$ws.Cells.Item(11, 5).Formula = '  -2.77%  '
$ws.Cells.Item(12, 4).Formula = '''0.497'
$ws.Cells.Item(12, 5).Formula = '  +6.23%  '
$ws.Cells.Item(13, 4).Formula = '''40.19'
$ws.Cells.Item(13, 5).Formula = '  +9.72%  '
$ws.Cells.Item(14, 4).Formula = '''0.0000255'
$ws.Cells.Item(14, 5).Formula = '  +5.37%  '
$ws.Cells.Item(15, 4).Formula = '4.406.29'
$ws.Cells.Item(15, 5).Formula = '  +21.16%  '
$ws.Cells.Item(16, 4).Formula = '3.779.39'
$ws.Cells.Item(16, 5).Formula = '  +21.39%  '
$ws.Cells.Item(17, 4).Formula = '70.073.49'
$ws.Cells.Item(17, 5).Formula = '  +4.04%  '
$ws.Cells.Item(18, 4).Formula = '''0.123'
$ws.Cells.Item(18, 5).Formula = '  +0.84%  '
$ws.Cells.Item(19, 4).Formula = '''7.54'
$ws.Cells.Item(19, 5).Formula = '  +6.98%  '
$ws.Cells.Item(20, 4).Formula = '''519.23'
$ws.Cells.Item(20, 5).Formula = '  +6.95%  '
$ws.Cells.Item(21, 4).Formula = '''16.61'
$ws.Cells.Item(21, 5).Formula = '  +0.78%  '
$ws.Cells.Item(22, 4).Formula = '''9.36'
$ws.Cells.Item(22, 5).Formula = '  +21.24%  '
$ws.Cells.Item(23, 4).Formula = '''0.740'
$ws.Cells.Item(23, 5).Formula = '  +7.27%  '
$ws.Cells.Item(24, 4).Formula = '''88.39'
$ws.Cells.Item(24, 5).Formula = '  +5.67%  '
$ws.Cells.Item(25, 4).Formula = '''2.47'
$ws.Cells.Item(25, 5).Formula = '  +7.78%  '
$ws.Cells.Item(26, 4).Formula = '''13.41'
$ws.Cells.Item(26, 5).Formula = '  +5.35%  '
$ws.Cells.Item(27, 4).Formula = '''10.85'
$ws.Cells.Item(27, 5).Formula = '  +2.70%  '
$ws.Cells.Item(28, 4).Formula = '''0.999'
$ws.Cells.Item(28, 5).Formula = '  -0.14%  '
$ws.Cells.Item(29, 4).Formula = '''0.0000123'
$ws.Cells.Item(29, 5).Formula = '  +30.18%  '
$ws.Cells.Item(30, 4).Formula = '''2.50'
$ws.Cells.Item(30, 5).Formula = '  +7.26%  '
$ws.Cells.Item(31, 4).Formula = '''2.85'
$ws.Cells.Item(31, 5).Formula = '  +8.90%  '
$ws.Cells.Item(32, 4).Formula = '''7.88'
$ws.Cells.Item(32, 5).Formula = '  -1.76%  '
$ws.Cells.Item(33, 4).Formula = '''31.79'
$ws.Cells.Item(33, 5).Formula = '  +13.00%  '
$ws.Cells.Item(34, 4).Formula = '''0.114'
$ws.Cells.Item(34, 5).Formula = '  +1.63%  '
$ws.Cells.Item(35, 5).Formula = '  +0.09%  '
$ws.Cells.Item(36, 2).Formula = 'Mantle'
$ws.Cells.Item(36, 3).Formula = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(36, 4).Formula = '''1.05'
$ws.Cells.Item(36, 5).Formula = '  +9.82%  '
$ws.Cells.Item(37, 2).Formula = 'Filecoin'
$ws.Cells.Item(37, 3).Formula = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(37, 4).Formula = '''6.15'
$ws.Cells.Item(37, 5).Formula = '  +9.87%  '
$ws.Cells.Item(38, 4).Formula = '''0.340'
$ws.Cells.Item(38, 5).Formula = '  +5.97%  '
$ws.Cells.Item(39, 2).Formula = 'Stacks'
$ws.Cells.Item(39, 3).Formula = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39, 4).Formula = '''2.16'
$ws.Cells.Item(39, 5).Formula = '  +6.96%  '
$ws.Cells.Item(40, 2).Formula = 'Kaspa'
$ws.Cells.Item(40, 3).Formula = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).Formula = '''0.132'
$ws.Cells.Item(40, 5).Formula = '  +7.32%  '
$ws.Cells.Item(41, 4).Formula = '''51.41'
$ws.Cells.Item(41, 5).Formula = '  +4.39%  '
$ws.Cells.Item(42, 2).Formula = 'Arweave'
$ws.Cells.Item(42, 3).Formula = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(42, 4).Formula = '''44.34'
$ws.Cells.Item(42, 5).Formula = '  -8.42%  '
$ws.Cells.Item(43, 2).Formula = 'Maker'
$ws.Cells.Item(43, 3).Formula = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Formula = '3.119.82'
$ws.Cells.Item(43, 5).Formula = '  +11.68%  '
$ws.Cells.Item(44, 2).Formula = 'Cosmos'
$ws.Cells.Item(44, 3).Formula = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(44, 4).Formula = '''8.79'
$ws.Cells.Item(44, 5).Formula = '  +5.66%  '
$ws.Cells.Item(45, 4).Formula = '''424.13'
$ws.Cells.Item(45, 5).Formula = '  +13.23%  '
$ws.Cells.Item(46, 4).Formula = '''2.71'
$ws.Cells.Item(46, 5).Formula = '  +0.30%  '
$ws.Cells.Item(47, 4).Formula = '''0.0365'
$ws.Cells.Item(47, 5).Formula = '  +5.37%  '
$ws.Cells.Item(48, 4).Formula = '''27.56'
$ws.Cells.Item(48, 5).Formula = '  +3.00%  '
$ws.Cells.Item(49, 4).Formula = '''137.18'
$ws.Cells.Item(49, 5).Formula = '  +1.06%  '
$ws.Cells.Item(50, 4).Formula = '''2.50'
$ws.Cells.Item(50, 5).Formula = '  +5.75%  '
